# Adds four new "example" rows (E0040-E0043) to the Example sheet,
# describing purpose-related examples, per the commit message:
# "adds examples for purposes" / "adds purpose module documentation".

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Example")

# Copy the formatting of the last existing data row (40) down into the
# four new rows (41-44) so they pick up the same cell styles. Column H
# is skipped (left blank, as on all other data rows).
$ws.Range("A40:G40").Copy()
$ws.Range("A41:G41").PasteSpecial(-4122) # xlPasteFormats
$ws.Range("A42:G42").PasteSpecial(-4122)
$ws.Range("A43:G43").PasteSpecial(-4122)
$ws.Range("A44:G44").PasteSpecial(-4122)

$ws.Range("I40:K40").Copy()
$ws.Range("I41:K41").PasteSpecial(-4122)
$ws.Range("I42:K42").PasteSpecial(-4122)
$ws.Range("I43:K43").PasteSpecial(-4122)
$ws.Range("I44:K44").PasteSpecial(-4122)

# Row 41 - E0040
$ws.Range("A41").Value = "E0040"
$ws.Range("B41").Value = "Extending a purpose and using human-readable descriptions"
$ws.Range("C41").Value = "This example describes how a purpose can be made clearer and accurate by two methods: (1) providing a human-readable description - which is always recommended; and (2) by extending a DPV concept"
$ws.Range("D41").Value = "E0040.ttl"
$ws.Range("E41").Value = "ttl"
$ws.Range("F41").Value = "file"
$ws.Range("G41").Value = "dpv:Purpose"
$ws.Range("I41").Value = "accepted"
$ws.Range("J41").Value = 45454.0
$ws.Range("K41").Value = "Harshvardhan J. Pandit"

# Row 42 - E0041
$ws.Range("A42").Value = "E0041"
$ws.Range("B42").Value = "Indicating purposes associated with a Service"
$ws.Range("C42").Value = "This example describes how the different purposes and information associated with a service can be expressed in a modular and clear manner"
$ws.Range("D42").Value = "E0041.ttl"
$ws.Range("E42").Value = "ttl"
$ws.Range("F42").Value = "file"
$ws.Range("G42").Value = "dpv:Purpose,dpv:Service,dpv:Process"
$ws.Range("I42").Value = "accepted"
$ws.Range("J42").Value = 45454.0
$ws.Range("K42").Value = "Harshvardhan J. Pandit"

# Row 43 - E0042
$ws.Range("A43").Value = "E0042"
$ws.Range("B43").Value = "Indicating legal compliance as a purpose along with the relevant law"
$ws.Range("C43").Value = "This example describes a purpose for performing 'KYC' identity verification as part of legal compliance with anti-money laundering laws"
$ws.Range("D43").Value = "E0042.ttl"
$ws.Range("E43").Value = "ttl"
$ws.Range("F43").Value = "file"
$ws.Range("G43").Value = "dpv:Purpose,dpv:LegalObligation"
$ws.Range("I43").Value = "accepted"
$ws.Range("J43").Value = 45454.0
$ws.Range("K43").Value = "Harshvardhan J. Pandit"

# Row 44 - E0043
$ws.Range("A44").Value = "E0043"
$ws.Range("B44").Value = "Indicating sector or domain and associating it with a purpose"
$ws.Range("C44").Value = "This example uses the NACE taxonomy to indicate the domain or sectorial relevance of the purpose"
$ws.Range("D44").Value = "E0043.ttl"
$ws.Range("E44").Value = "ttl"
$ws.Range("F44").Value = "file"
$ws.Range("G44").Value = "dpv:Purpose,dpv:Sector"
$ws.Range("I44").Value = "accepted"
$ws.Range("J44").Value = 45454.0
$ws.Range("K44").Value = "Harshvardhan J. Pandit"
